$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 81

# Copy formatting (style) from the row above for the styled columns (A and E)
$ws.Range("A80").Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)
$ws.Range("E80").Copy()
$ws.Range("E" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A" + $newRow).Value = 80
$ws.Range("B" + $newRow).Value = "bosnia-and-herzegovina"
$ws.Range("C" + $newRow).Value = "premijer-liga-bih"
$ws.Range("D" + $newRow).Value = "2023-2024"
$ws.Range("E" + $newRow).Value = 45242.63541666666
$ws.Range("F" + $newRow).Value = "Borac Banja Luka"
$ws.Range("G" + $newRow).Value = 3
$ws.Range("H" + $newRow).Value = "Zeljeznicar"
$ws.Range("I" + $newRow).Value = 2
$ws.Range("J" + $newRow).Value = 1.52
$ws.Range("K" + $newRow).Value = "12/11/2023 05:12"
$ws.Range("L" + $newRow).Value = 1.43
$ws.Range("M" + $newRow).Value = "12/11/2023 15:09"
$ws.Range("N" + $newRow).Value = 3.79
$ws.Range("O" + $newRow).Value = "12/11/2023 05:12"
$ws.Range("P" + $newRow).Value = 4.25
$ws.Range("Q" + $newRow).Value = "12/11/2023 15:09"
$ws.Range("R" + $newRow).Value = 5.81
$ws.Range("S" + $newRow).Value = "12/11/2023 05:12"
$ws.Range("T" + $newRow).Value = 7.46
$ws.Range("U" + $newRow).Value = "12/11/2023 15:09"
$ws.Range("V" + $newRow).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/borac-banja-luka-zeljeznicar/hOMeLggh/"
